$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "model_res" column (I) mirroring the existing "logic_res" column (G/H)
$ws.Range("I1").Value = "model_res"
$ws.Range("I2").Value = "string"
$ws.Range("I3").Value = "model逻辑"

# Match the phonetic-range flag used by the other header/meta cells in the sheet
$ws.Range("I1:I3").SetPhonetic()

# Update selection to the full column I, matching the committed workbook view
$ws.Columns.Item(9).Select()
